$wb = $excel.ActiveWorkbook

# --- Shared string text tweaks (append an extra 'z') ---
$zs44 = ""
for ($i = 0; $i -lt 44; $i++) {
  $zs44 = $zs44 + "z"
}
$zs84 = ""
for ($i = 0; $i -lt 84; $i++) {
  $zs84 = $zs84 + "z"
}

$wsBudgetOut = $wb.Worksheets.Item("Budget Out")
$wsBudgetOut.Range("F9").Value = "Description007" + $zs44

$wsTestRecord = $wb.Worksheets.Item("TestRecord")
$wsTestRecord.Range("E10").Value = "some test text" + $zs84

# --- Numeric / date updates ---

# TestRecord sheet: row 10
$wsTestRecord.Range("A10").Value = 43263
$wsTestRecord.Range("B10").Value = 121.14

# Budget Out sheet: row 9
$wsBudgetOut.Range("C9").Value = 90.22

# Expected Out sheet: rows 9 and 11 (B1 SUM formula recalculates automatically)
$wsExpectedOut = $wb.Worksheets.Item("Expected Out")
$wsExpectedOut.Range("B9").Value = 1348.16
$wsExpectedOut.Range("B11").Value = 428.02

$excel.Calculate()
